# Weekly fruit/vegetable price update: a new price-report row for
# "Femacal de La Calera" / Ajo (Coquimbo) is inserted at row 213,
# pushing the existing rows 213-245 down to 214-246.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 213 (shifts old rows 213..245 -> 214..246)
$ws.Rows.Item(213).Insert()

# Populate the newly inserted row 213 with the new week's data
$ws.Range("A213").Value2 = 3
$ws.Range("B213").Value2 = "Femacal de La Calera"
$ws.Range("C213").Value2 = "Coquimbo"
$ws.Range("D213").Value2 = 44505
$ws.Range("D213").NumberFormat = $ws.Range("D214").NumberFormat
$ws.Range("E213").Value2 = 5
$ws.Range("F213").Value2 = 100112003
$ws.Range("G213").Value2 = "Ajo"
$ws.Range("H213").Value2 = "Chino"
$ws.Range("I213").Value2 = "Primera"
$ws.Range("J213").Value2 = 65
$ws.Range("K213").Value2 = 16000
$ws.Range("L213").Value2 = 16500
$ws.Range("M213").Value2 = 16231
$ws.Range("N213").Value2 = "$/caja 10 kilos"
$ws.Range("O213").Value2 = "China"
$ws.Range("P213").Value2 = 1623
$ws.Range("Q213").Value2 = 10
$ws.Range("R213").Value2 = "Hortaliza"
